$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix for the naive component forecaster: each data row (rows 2-16, columns
# B:K) was missing its most-recent error value. For every row we now insert a
# new "first" value into column B and shift the row's existing values one
# column to the right. Because the sheet's populated region stops at column K,
# any shifted value that would land past column K is simply dropped (this only
# affects the fully-populated rows 2-6, which lose their old column-K value).
$newLeadingValues = @{
    2  = 0.3648791949059138
    3  = -0.2352699264540507
    4  = -0.05148746350304451
    5  = -0.1333319740152609
    6  = 1.614150253737389
    7  = 0.5701030647716323
    8  = 0.2202779152847414
    9  = 0.5040960054549828
    10 = 0.420735823599318
    11 = -0.1252583916527783
    12 = 0.08824118641116785
    13 = -0.1133200159455487
    14 = 0.1743923273248104
    15 = -0.4559694969238889
    16 = 0.1808172637304477
}

$firstCol = 2   # column B
$maxCol   = 11  # column K is the sheet's existing right edge

for ($row = 2; $row -le 16; $row++) {

    # Find the last populated column in this row within B:K.
    $lastCol = $firstCol - 1
    for ($col = $maxCol; $col -ge $firstCol; $col--) {
        $cellVal = $ws.Cells.Item($row, $col).Value()
        if ($cellVal -ne $null) {
            $lastCol = $col
            break
        }
    }

    # Shift existing values one column to the right, starting from the
    # rightmost populated cell so we never overwrite a value before reading it.
    for ($col = $lastCol; $col -ge $firstCol; $col--) {
        $srcVal  = $ws.Cells.Item($row, $col).Value()
        $destCol = $col + 1
        if ($destCol -le $maxCol) {
            $ws.Cells.Item($row, $destCol).Value = $srcVal
        }
    }

    # Place the new leading value in column B.
    $ws.Cells.Item($row, $firstCol).Value = $newLeadingValues[$row]
}
